$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.848.40'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -4.23%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.436.78'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -5.39%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.22'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.82'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -6.86%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.598'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -5.88%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.421.93'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -5.54%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.181'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -7.40%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.65'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.85%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.563'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -9.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '45.89'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -5.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000266'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -5.97%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.004.24'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -5.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.19'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -9.10%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '602.14'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -9.99%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.119.46'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.455.48'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -4.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.119'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.96'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.86%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.78'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -6.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.862'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -8.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.16'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -11.81%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '94.44'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -5.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.74'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.24%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.54'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -9.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.87'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -9.90%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.76'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -8.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.27'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -9.76%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.01'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -8.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.28'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -7.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.72'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -12.14%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '610.27'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +6.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.43'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -5.92%  '
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'FirstDigitalUSD'
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.62%  '
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.97'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.26%  '
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'dogwifhat'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.35'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -16.00%  '
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0989'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -8.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0425'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -5.92%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.132'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -5.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.319.56'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -6.48%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.317'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -8.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '31.85'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -7.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₃0677'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -7.27%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.44'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -8.68%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.66'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -13.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.126'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -6.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '132.94'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.38%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'USDe'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.01%  '
